$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for rows 2-5 (A and B columns)
$rows = 2..5
$data = @()
foreach ($r in $rows) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $data += ,@($a, $b)
}

# Rotate: row2 data moves to the end (row5), rows 3,4,5 shift up to 2,3,4
$rotated = @()
for ($i = 1; $i -lt $data.Count; $i++) {
    $rotated += ,$data[$i]
}
$rotated += ,$data[0]

# Write back rotated values
for ($i = 0; $i -lt $rotated.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rotated[$i][0]
    $ws.Cells.Item($r, 2).Value = $rotated[$i][1]
}
